# Append two new staff rows (Customer_Service role) to Sheet1, each with a
# date-of-birth formatted as a short date, and move the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: Do Xuan Thanh (Male)
$ws.Range("A23").Value = "9edbd47a-186b-4c5b-a077-8446b7418f6f"
$ws.Range("B23").Value = "xuanthanh"
$ws.Range("C23").Value = "Customer_Service"
$ws.Range("D23").Value = "xuanthanh@gmail.com"
$ws.Range("E23").Value = 346494851
$ws.Range("F23").Value = '$2y$10$Sgk31oRzu3kNLfIUvAJLNOaTwts89qxic3Yzk2s59o0C/yMxabh2G'
$ws.Range("G23").Value = "Đỗ Xuân Thanh"
$ws.Range("H23").Value = "https://vapa.vn/wp-content/uploads/2022/12/anh-avatar-facebook-dep-001.jpg"
$ws.Range("I23").Value = "Male"
$ws.Range("J23").Value = 37011
$ws.Range("J23").NumberFormat = "mm-dd-yy"
$ws.Range("K23").Value = "Bình Dương"
$ws.Range("L23").Value = $true
$ws.Range("M23").Value = $true
$ws.Range("N23").Value = "Đại học Y khoa Vinh"

# Row 24: Le Ngoc Nhu (Female)
$ws.Range("A24").Value = "2798c948-07a7-4f85-b7a6-8d8d69e53676"
$ws.Range("B24").Value = "ngocnhu"
$ws.Range("C24").Value = "Customer_Service"
$ws.Range("D24").Value = "ngocnhu@gmail.com"
$ws.Range("E24").Value = 366995813
$ws.Range("F24").Value = '$2y$10$Sgk31oRzu3kNLfIUvAJLNOaTwts89qxic3Yzk2s59o0C/yMxabh2G'
$ws.Range("G24").Value = "Lê Ngọc Như"
$ws.Range("H24").Value = "https://thao68.com/wp-content/uploads/2022/03/avatar-facebook-3.jpg"
$ws.Range("I24").Value = "Female"
$ws.Range("J24").Value = 37001
$ws.Range("J23").Copy()
$ws.Range("J24").PasteSpecial(-4122)
$ws.Range("K24").Value = "Q8, TP.HCM"
$ws.Range("L24").Value = $true
$ws.Range("M24").Value = $true
$ws.Range("N24").Value = "Đại học Y khoa Vinh"

# Move the selection to reflect where the user last clicked after the edit.
$ws.Range("J25").Select() | Out-Null
